$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 232
$ws.Range("F4").Value = 69
$ws.Range("F5").Value = 9226
$ws.Range("F6").Value = 9226
$ws.Range("F7").Value = 560
$ws.Range("F9").Value = 152
$ws.Range("F10").Value = 233
$ws.Range("F12").Value = 388
$ws.Range("F13").Value = 145
$ws.Range("F14").Value = 150
$ws.Range("F15").Value = 419
$ws.Range("F16").Value = 11795
$ws.Range("F17").Value = 11795
$ws.Range("F27").Value = 169
$ws.Range("F29").Value = 2711
$ws.Range("F33").Value = 60
$ws.Range("F35").Value = 2137
$ws.Range("F36").Value = 972
$ws.Range("F37").Value = 4175
$ws.Range("F39").Value = 3595
$ws.Range("F40").Value = 340
$ws.Range("F41").Value = 2609
$ws.Range("F43").Value = 1300
$ws.Range("F44").Value = 186
$ws.Range("F45").Value = 767
$ws.Range("F46").Value = 402
$ws.Range("F47").Value = 465
$ws.Range("F48").Value = 61
$ws.Range("F49").Value = 196
$ws.Range("F50").Value = 118
$ws.Range("F51").Value = 115

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value = 47
$ws.Range("F15").Value = 42
$ws.Range("F17").Value = 3
$ws.Range("F21").Value = 73

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 45

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F9").Value = 232
$ws.Range("F10").Value = 9226
$ws.Range("F11").Value = 9226
$ws.Range("F12").Value = 560
$ws.Range("F14").Value = 152
$ws.Range("F15").Value = 233
$ws.Range("F16").Value = 388
$ws.Range("F17").Value = 145
$ws.Range("F18").Value = 150
$ws.Range("F19").Value = 11795
$ws.Range("F20").Value = 11795
$ws.Range("F24").Value = 45
$ws.Range("F30").Value = 169
$ws.Range("F32").Value = 2711
$ws.Range("F36").Value = 60
$ws.Range("F38").Value = 2137
$ws.Range("F39").Value = 972
$ws.Range("F42").Value = 3595
$ws.Range("F44").Value = 73
$ws.Range("F45").Value = 1300
$ws.Range("F46").Value = 186
$ws.Range("F47").Value = 402
$ws.Range("F49").Value = 465
$ws.Range("F50").Value = 61
$ws.Range("F51").Value = 196

